$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-93: columns D (date), J (volume), K/L/M (prices), N (unit), P (price/kg)
$ws.Cells.Item(2, 4).Value = 44706
$ws.Cells.Item(2, 10).Value = 500
$ws.Cells.Item(2, 11).Value = 13000
$ws.Cells.Item(2, 12).Value = 14000
$ws.Cells.Item(2, 13).Value = 13500
$ws.Cells.Item(2, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(2, 16).Value = 540
$ws.Cells.Item(3, 4).Value = 44707
$ws.Cells.Item(3, 10).Value = 500
$ws.Cells.Item(3, 11).Value = 13000
$ws.Cells.Item(3, 12).Value = 14000
$ws.Cells.Item(3, 13).Value = 13500
$ws.Cells.Item(3, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(3, 16).Value = 540
$ws.Cells.Item(4, 4).Value = 44329
$ws.Cells.Item(4, 10).Value = 1000
$ws.Cells.Item(4, 11).Value = 12000
$ws.Cells.Item(4, 12).Value = 13000
$ws.Cells.Item(4, 13).Value = 12500
$ws.Cells.Item(4, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(4, 16).Value = 500
$ws.Cells.Item(5, 4).Value = 44447
$ws.Cells.Item(5, 10).Value = 1000
$ws.Cells.Item(5, 11).Value = 10000
$ws.Cells.Item(5, 12).Value = 12000
$ws.Cells.Item(5, 13).Value = 11000
$ws.Cells.Item(5, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(5, 16).Value = 440
$ws.Cells.Item(6, 4).Value = 44510
$ws.Cells.Item(6, 10).Value = 1300
$ws.Cells.Item(6, 11).Value = 6000
$ws.Cells.Item(6, 12).Value = 7000
$ws.Cells.Item(6, 13).Value = 6500
$ws.Cells.Item(6, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(6, 16).Value = 260
$ws.Cells.Item(7, 4).Value = 44756
$ws.Cells.Item(7, 10).Value = 500
$ws.Cells.Item(7, 11).Value = 11000
$ws.Cells.Item(7, 12).Value = 12000
$ws.Cells.Item(7, 13).Value = 11500
$ws.Cells.Item(7, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(7, 16).Value = 460
$ws.Cells.Item(8, 4).Value = 44482
$ws.Cells.Item(8, 10).Value = 1600
$ws.Cells.Item(8, 11).Value = 4000
$ws.Cells.Item(8, 12).Value = 5000
$ws.Cells.Item(8, 13).Value = 4500
$ws.Cells.Item(8, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(8, 16).Value = 180
$ws.Cells.Item(9, 4).Value = 44847
$ws.Cells.Item(9, 10).Value = 800
$ws.Cells.Item(9, 11).Value = 4500
$ws.Cells.Item(9, 12).Value = 5000
$ws.Cells.Item(9, 13).Value = 4750
$ws.Cells.Item(9, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(9, 16).Value = 190
$ws.Cells.Item(10, 4).Value = 44839
$ws.Cells.Item(10, 10).Value = 1200
$ws.Cells.Item(10, 11).Value = 5000
$ws.Cells.Item(10, 12).Value = 6000
$ws.Cells.Item(10, 13).Value = 5500
$ws.Cells.Item(10, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(10, 16).Value = 220
$ws.Cells.Item(11, 4).Value = 44419
$ws.Cells.Item(11, 10).Value = 1100
$ws.Cells.Item(11, 11).Value = 11000
$ws.Cells.Item(11, 12).Value = 12000
$ws.Cells.Item(11, 13).Value = 11500
$ws.Cells.Item(11, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(11, 16).Value = 460
$ws.Cells.Item(12, 4).Value = 44462
$ws.Cells.Item(12, 10).Value = 800
$ws.Cells.Item(12, 11).Value = 9000
$ws.Cells.Item(12, 12).Value = 10000
$ws.Cells.Item(12, 13).Value = 9500
$ws.Cells.Item(12, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(12, 16).Value = 380
$ws.Cells.Item(13, 4).Value = 44798
$ws.Cells.Item(13, 10).Value = 1500
$ws.Cells.Item(13, 11).Value = 7500
$ws.Cells.Item(13, 12).Value = 8000
$ws.Cells.Item(13, 13).Value = 7750
$ws.Cells.Item(13, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(13, 16).Value = 310
$ws.Cells.Item(14, 4).Value = 44356
$ws.Cells.Item(14, 10).Value = 1000
$ws.Cells.Item(14, 11).Value = 11000
$ws.Cells.Item(14, 12).Value = 12000
$ws.Cells.Item(14, 13).Value = 11500
$ws.Cells.Item(14, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(14, 16).Value = 460
$ws.Cells.Item(15, 4).Value = 44783
$ws.Cells.Item(15, 10).Value = 1800
$ws.Cells.Item(15, 11).Value = 8000
$ws.Cells.Item(15, 12).Value = 9000
$ws.Cells.Item(15, 13).Value = 8500
$ws.Cells.Item(15, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(15, 16).Value = 340
$ws.Cells.Item(16, 4).Value = 44363
$ws.Cells.Item(16, 10).Value = 900
$ws.Cells.Item(16, 11).Value = 11000
$ws.Cells.Item(16, 12).Value = 12000
$ws.Cells.Item(16, 13).Value = 11500
$ws.Cells.Item(16, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(16, 16).Value = 460
$ws.Cells.Item(17, 4).Value = 44818
$ws.Cells.Item(17, 10).Value = 1100
$ws.Cells.Item(17, 11).Value = 7000
$ws.Cells.Item(17, 12).Value = 8000
$ws.Cells.Item(17, 13).Value = 7500
$ws.Cells.Item(17, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(17, 16).Value = 300
$ws.Cells.Item(18, 4).Value = 44804
$ws.Cells.Item(18, 10).Value = 1000
$ws.Cells.Item(18, 11).Value = 8000
$ws.Cells.Item(18, 12).Value = 9000
$ws.Cells.Item(18, 13).Value = 8500
$ws.Cells.Item(18, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(18, 16).Value = 340
$ws.Cells.Item(19, 4).Value = 44721
$ws.Cells.Item(19, 10).Value = 500
$ws.Cells.Item(19, 11).Value = 12000
$ws.Cells.Item(19, 12).Value = 14000
$ws.Cells.Item(19, 13).Value = 13000
$ws.Cells.Item(19, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(19, 16).Value = 520
$ws.Cells.Item(20, 4).Value = 44328
$ws.Cells.Item(20, 10).Value = 900
$ws.Cells.Item(20, 11).Value = 11000
$ws.Cells.Item(20, 12).Value = 12000
$ws.Cells.Item(20, 13).Value = 11500
$ws.Cells.Item(20, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(20, 16).Value = 460
$ws.Cells.Item(21, 4).Value = 44819
$ws.Cells.Item(21, 10).Value = 1000
$ws.Cells.Item(21, 11).Value = 8000
$ws.Cells.Item(21, 12).Value = 9000
$ws.Cells.Item(21, 13).Value = 8500
$ws.Cells.Item(21, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(21, 16).Value = 340
$ws.Cells.Item(22, 4).Value = 44343
$ws.Cells.Item(22, 10).Value = 500
$ws.Cells.Item(22, 11).Value = 9000
$ws.Cells.Item(22, 12).Value = 10000
$ws.Cells.Item(22, 13).Value = 9500
$ws.Cells.Item(22, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(22, 16).Value = 380
$ws.Cells.Item(23, 4).Value = 44426
$ws.Cells.Item(23, 10).Value = 500
$ws.Cells.Item(23, 11).Value = 11000
$ws.Cells.Item(23, 12).Value = 12000
$ws.Cells.Item(23, 13).Value = 11500
$ws.Cells.Item(23, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(23, 16).Value = 460
$ws.Cells.Item(24, 4).Value = 44336
$ws.Cells.Item(24, 10).Value = 1200
$ws.Cells.Item(24, 11).Value = 12000
$ws.Cells.Item(24, 12).Value = 13000
$ws.Cells.Item(24, 13).Value = 12500
$ws.Cells.Item(24, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(24, 16).Value = 500
$ws.Cells.Item(25, 4).Value = 44889
$ws.Cells.Item(25, 10).Value = 520
$ws.Cells.Item(25, 11).Value = 5000
$ws.Cells.Item(25, 12).Value = 6000
$ws.Cells.Item(25, 13).Value = 5500
$ws.Cells.Item(25, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(25, 16).Value = 220
$ws.Cells.Item(26, 4).Value = 44469
$ws.Cells.Item(26, 10).Value = 600
$ws.Cells.Item(26, 11).Value = 5000
$ws.Cells.Item(26, 12).Value = 6000
$ws.Cells.Item(26, 13).Value = 5500
$ws.Cells.Item(26, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(26, 16).Value = 220
$ws.Cells.Item(27, 4).Value = 44385
$ws.Cells.Item(27, 10).Value = 600
$ws.Cells.Item(27, 11).Value = 8000
$ws.Cells.Item(27, 12).Value = 9000
$ws.Cells.Item(27, 13).Value = 8500
$ws.Cells.Item(27, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(27, 16).Value = 340
$ws.Cells.Item(28, 4).Value = 44489
$ws.Cells.Item(28, 10).Value = 1200
$ws.Cells.Item(28, 11).Value = 5000
$ws.Cells.Item(28, 12).Value = 6000
$ws.Cells.Item(28, 13).Value = 5500
$ws.Cells.Item(28, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(28, 16).Value = 220
$ws.Cells.Item(29, 4).Value = 44455
$ws.Cells.Item(29, 10).Value = 600
$ws.Cells.Item(29, 11).Value = 9000
$ws.Cells.Item(29, 12).Value = 10000
$ws.Cells.Item(29, 13).Value = 9500
$ws.Cells.Item(29, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(29, 16).Value = 380
$ws.Cells.Item(30, 4).Value = 45077
$ws.Cells.Item(30, 10).Value = 700
$ws.Cells.Item(30, 11).Value = 12000
$ws.Cells.Item(30, 12).Value = 14000
$ws.Cells.Item(30, 13).Value = 13000
$ws.Cells.Item(30, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(30, 16).Value = 520
$ws.Cells.Item(31, 4).Value = 44399
$ws.Cells.Item(31, 10).Value = 500
$ws.Cells.Item(31, 11).Value = 9000
$ws.Cells.Item(31, 12).Value = 10000
$ws.Cells.Item(31, 13).Value = 9500
$ws.Cells.Item(31, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(31, 16).Value = 380
$ws.Cells.Item(32, 4).Value = 44476
$ws.Cells.Item(32, 10).Value = 1100
$ws.Cells.Item(32, 11).Value = 5000
$ws.Cells.Item(32, 12).Value = 6000
$ws.Cells.Item(32, 13).Value = 5500
$ws.Cells.Item(32, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(32, 16).Value = 220
$ws.Cells.Item(33, 4).Value = 44349
$ws.Cells.Item(33, 10).Value = 600
$ws.Cells.Item(33, 11).Value = 10000
$ws.Cells.Item(33, 12).Value = 12000
$ws.Cells.Item(33, 13).Value = 11000
$ws.Cells.Item(33, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(33, 16).Value = 440
$ws.Cells.Item(34, 4).Value = 45063
$ws.Cells.Item(34, 10).Value = 700
$ws.Cells.Item(34, 11).Value = 12000
$ws.Cells.Item(34, 12).Value = 14000
$ws.Cells.Item(34, 13).Value = 13000
$ws.Cells.Item(34, 14).Value = '$/malla 25 kilos'
$ws.Cells.Item(34, 16).Value = 520
$ws.Cells.Item(35, 4).Value = 44776
$ws.Cells.Item(35, 10).Value = 1100
$ws.Cells.Item(35, 11).Value = 10000
$ws.Cells.Item(35, 12).Value = 11000
$ws.Cells.Item(35, 13).Value = 10500
$ws.Cells.Item(35, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(35, 16).Value = 420
$ws.Cells.Item(36, 4).Value = 44868
$ws.Cells.Item(36, 10).Value = 1100
$ws.Cells.Item(36, 11).Value = 4000
$ws.Cells.Item(36, 12).Value = 5000
$ws.Cells.Item(36, 13).Value = 4500
$ws.Cells.Item(36, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(36, 16).Value = 180
$ws.Cells.Item(37, 4).Value = 44832
$ws.Cells.Item(37, 10).Value = 700
$ws.Cells.Item(37, 11).Value = 6000
$ws.Cells.Item(37, 12).Value = 7000
$ws.Cells.Item(37, 13).Value = 6500
$ws.Cells.Item(37, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(37, 16).Value = 260
$ws.Cells.Item(38, 4).Value = 44461
$ws.Cells.Item(38, 10).Value = 1100
$ws.Cells.Item(38, 11).Value = 9000
$ws.Cells.Item(38, 12).Value = 10000
$ws.Cells.Item(38, 13).Value = 9500
$ws.Cells.Item(38, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(38, 16).Value = 380
$ws.Cells.Item(39, 4).Value = 44377
$ws.Cells.Item(39, 10).Value = 800
$ws.Cells.Item(39, 11).Value = 9000
$ws.Cells.Item(39, 12).Value = 10000
$ws.Cells.Item(39, 13).Value = 9500
$ws.Cells.Item(39, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(39, 16).Value = 380
$ws.Cells.Item(40, 4).Value = 45084
$ws.Cells.Item(40, 10).Value = 1100
$ws.Cells.Item(40, 11).Value = 11000
$ws.Cells.Item(40, 12).Value = 13000
$ws.Cells.Item(40, 13).Value = 12000
$ws.Cells.Item(40, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(40, 16).Value = 480
$ws.Cells.Item(41, 4).Value = 45091
$ws.Cells.Item(41, 10).Value = 1100
$ws.Cells.Item(41, 11).Value = 13000
$ws.Cells.Item(41, 12).Value = 14000
$ws.Cells.Item(41, 13).Value = 13500
$ws.Cells.Item(41, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(41, 16).Value = 540
$ws.Cells.Item(42, 4).Value = 44454
$ws.Cells.Item(42, 10).Value = 800
$ws.Cells.Item(42, 11).Value = 9000
$ws.Cells.Item(42, 12).Value = 10000
$ws.Cells.Item(42, 13).Value = 9500
$ws.Cells.Item(42, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(42, 16).Value = 380
$ws.Cells.Item(43, 4).Value = 44371
$ws.Cells.Item(43, 10).Value = 500
$ws.Cells.Item(43, 11).Value = 10000
$ws.Cells.Item(43, 12).Value = 12000
$ws.Cells.Item(43, 13).Value = 11000
$ws.Cells.Item(43, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(43, 16).Value = 440
$ws.Cells.Item(44, 4).Value = 44475
$ws.Cells.Item(44, 10).Value = 1200
$ws.Cells.Item(44, 11).Value = 5000
$ws.Cells.Item(44, 12).Value = 6000
$ws.Cells.Item(44, 13).Value = 5500
$ws.Cells.Item(44, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(44, 16).Value = 220
$ws.Cells.Item(45, 4).Value = 44860
$ws.Cells.Item(45, 10).Value = 700
$ws.Cells.Item(45, 11).Value = 4000
$ws.Cells.Item(45, 12).Value = 5000
$ws.Cells.Item(45, 13).Value = 4500
$ws.Cells.Item(45, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(45, 16).Value = 180
$ws.Cells.Item(46, 4).Value = 44308
$ws.Cells.Item(46, 10).Value = 400
$ws.Cells.Item(46, 11).Value = 11000
$ws.Cells.Item(46, 12).Value = 12000
$ws.Cells.Item(46, 13).Value = 11500
$ws.Cells.Item(46, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(46, 16).Value = 460
$ws.Cells.Item(47, 4).Value = 44398
$ws.Cells.Item(47, 10).Value = 400
$ws.Cells.Item(47, 11).Value = 9000
$ws.Cells.Item(47, 12).Value = 10000
$ws.Cells.Item(47, 13).Value = 9500
$ws.Cells.Item(47, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(47, 16).Value = 380
$ws.Cells.Item(48, 4).Value = 45092
$ws.Cells.Item(48, 10).Value = 1200
$ws.Cells.Item(48, 11).Value = 13000
$ws.Cells.Item(48, 12).Value = 14000
$ws.Cells.Item(48, 13).Value = 13500
$ws.Cells.Item(48, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(48, 16).Value = 540
$ws.Cells.Item(49, 4).Value = 44910
$ws.Cells.Item(49, 10).Value = 500
$ws.Cells.Item(49, 11).Value = 7000
$ws.Cells.Item(49, 12).Value = 8000
$ws.Cells.Item(49, 13).Value = 7500
$ws.Cells.Item(49, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(49, 16).Value = 300
$ws.Cells.Item(50, 4).Value = 44448
$ws.Cells.Item(50, 10).Value = 800
$ws.Cells.Item(50, 11).Value = 10000
$ws.Cells.Item(50, 12).Value = 12000
$ws.Cells.Item(50, 13).Value = 11000
$ws.Cells.Item(50, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(50, 16).Value = 440
$ws.Cells.Item(51, 4).Value = 44846
$ws.Cells.Item(51, 10).Value = 1600
$ws.Cells.Item(51, 11).Value = 4500
$ws.Cells.Item(51, 12).Value = 5000
$ws.Cells.Item(51, 13).Value = 4750
$ws.Cells.Item(51, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(51, 16).Value = 190
$ws.Cells.Item(52, 4).Value = 44762
$ws.Cells.Item(52, 10).Value = 1500
$ws.Cells.Item(52, 11).Value = 11000
$ws.Cells.Item(52, 12).Value = 12000
$ws.Cells.Item(52, 13).Value = 11500
$ws.Cells.Item(52, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(52, 16).Value = 460
$ws.Cells.Item(53, 4).Value = 44882
$ws.Cells.Item(53, 10).Value = 560
$ws.Cells.Item(53, 11).Value = 4500
$ws.Cells.Item(53, 12).Value = 5000
$ws.Cells.Item(53, 13).Value = 4750
$ws.Cells.Item(53, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(53, 16).Value = 190
$ws.Cells.Item(54, 4).Value = 45113
$ws.Cells.Item(54, 10).Value = 1000
$ws.Cells.Item(54, 11).Value = 10000
$ws.Cells.Item(54, 12).Value = 12000
$ws.Cells.Item(54, 13).Value = 11000
$ws.Cells.Item(54, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(54, 16).Value = 440
$ws.Cells.Item(55, 4).Value = 44406
$ws.Cells.Item(55, 10).Value = 800
$ws.Cells.Item(55, 11).Value = 10000
$ws.Cells.Item(55, 12).Value = 11000
$ws.Cells.Item(55, 13).Value = 10500
$ws.Cells.Item(55, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(55, 16).Value = 420
$ws.Cells.Item(56, 4).Value = 44392
$ws.Cells.Item(56, 10).Value = 600
$ws.Cells.Item(56, 11).Value = 9000
$ws.Cells.Item(56, 12).Value = 10000
$ws.Cells.Item(56, 13).Value = 9500
$ws.Cells.Item(56, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(56, 16).Value = 380
$ws.Cells.Item(57, 4).Value = 44755
$ws.Cells.Item(57, 10).Value = 1100
$ws.Cells.Item(57, 11).Value = 11000
$ws.Cells.Item(57, 12).Value = 12000
$ws.Cells.Item(57, 13).Value = 11500
$ws.Cells.Item(57, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(57, 16).Value = 460
$ws.Cells.Item(58, 4).Value = 44335
$ws.Cells.Item(58, 10).Value = 1000
$ws.Cells.Item(58, 11).Value = 12000
$ws.Cells.Item(58, 12).Value = 13000
$ws.Cells.Item(58, 13).Value = 12500
$ws.Cells.Item(58, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(58, 16).Value = 500
$ws.Cells.Item(59, 4).Value = 44434
$ws.Cells.Item(59, 10).Value = 600
$ws.Cells.Item(59, 11).Value = 10000
$ws.Cells.Item(59, 12).Value = 11000
$ws.Cells.Item(59, 13).Value = 10500
$ws.Cells.Item(59, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(59, 16).Value = 420
$ws.Cells.Item(60, 4).Value = 45050
$ws.Cells.Item(60, 10).Value = 300
$ws.Cells.Item(60, 11).Value = 12000
$ws.Cells.Item(60, 12).Value = 13000
$ws.Cells.Item(60, 13).Value = 12500
$ws.Cells.Item(60, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(60, 16).Value = 500
$ws.Cells.Item(61, 4).Value = 44483
$ws.Cells.Item(61, 10).Value = 1200
$ws.Cells.Item(61, 11).Value = 4000
$ws.Cells.Item(61, 12).Value = 5000
$ws.Cells.Item(61, 13).Value = 4500
$ws.Cells.Item(61, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(61, 16).Value = 180
$ws.Cells.Item(62, 4).Value = 44435
$ws.Cells.Item(62, 10).Value = 600
$ws.Cells.Item(62, 11).Value = 10000
$ws.Cells.Item(62, 12).Value = 11000
$ws.Cells.Item(62, 13).Value = 10500
$ws.Cells.Item(62, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(62, 16).Value = 420
$ws.Cells.Item(63, 4).Value = 44441
$ws.Cells.Item(63, 10).Value = 1100
$ws.Cells.Item(63, 11).Value = 11000
$ws.Cells.Item(63, 12).Value = 12000
$ws.Cells.Item(63, 13).Value = 11500
$ws.Cells.Item(63, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(63, 16).Value = 460
$ws.Cells.Item(64, 4).Value = 44412
$ws.Cells.Item(64, 10).Value = 1000
$ws.Cells.Item(64, 11).Value = 10000
$ws.Cells.Item(64, 12).Value = 11000
$ws.Cells.Item(64, 13).Value = 10500
$ws.Cells.Item(64, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(64, 16).Value = 420
$ws.Cells.Item(65, 4).Value = 44812
$ws.Cells.Item(65, 10).Value = 600
$ws.Cells.Item(65, 11).Value = 5000
$ws.Cells.Item(65, 12).Value = 6000
$ws.Cells.Item(65, 13).Value = 5500
$ws.Cells.Item(65, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(65, 16).Value = 220
$ws.Cells.Item(66, 4).Value = 44427
$ws.Cells.Item(66, 10).Value = 360
$ws.Cells.Item(66, 11).Value = 10000
$ws.Cells.Item(66, 12).Value = 11000
$ws.Cells.Item(66, 13).Value = 10500
$ws.Cells.Item(66, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(66, 16).Value = 420
$ws.Cells.Item(67, 4).Value = 44503
$ws.Cells.Item(67, 10).Value = 760
$ws.Cells.Item(67, 11).Value = 5000
$ws.Cells.Item(67, 12).Value = 6000
$ws.Cells.Item(67, 13).Value = 5500
$ws.Cells.Item(67, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(67, 16).Value = 220
$ws.Cells.Item(68, 4).Value = 44769
$ws.Cells.Item(68, 10).Value = 1300
$ws.Cells.Item(68, 11).Value = 7000
$ws.Cells.Item(68, 12).Value = 8000
$ws.Cells.Item(68, 13).Value = 7500
$ws.Cells.Item(68, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(68, 16).Value = 300
$ws.Cells.Item(69, 4).Value = 44384
$ws.Cells.Item(69, 10).Value = 700
$ws.Cells.Item(69, 11).Value = 8000
$ws.Cells.Item(69, 12).Value = 9000
$ws.Cells.Item(69, 13).Value = 8500
$ws.Cells.Item(69, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(69, 16).Value = 340
$ws.Cells.Item(70, 4).Value = 44727
$ws.Cells.Item(70, 10).Value = 408
$ws.Cells.Item(70, 11).Value = 14000
$ws.Cells.Item(70, 12).Value = 15000
$ws.Cells.Item(70, 13).Value = 14510
$ws.Cells.Item(70, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(70, 16).Value = 580
$ws.Cells.Item(71, 4).Value = 44413
$ws.Cells.Item(71, 10).Value = 1200
$ws.Cells.Item(71, 11).Value = 10000
$ws.Cells.Item(71, 12).Value = 11000
$ws.Cells.Item(71, 13).Value = 10500
$ws.Cells.Item(71, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(71, 16).Value = 420
$ws.Cells.Item(72, 4).Value = 44699
$ws.Cells.Item(72, 10).Value = 400
$ws.Cells.Item(72, 11).Value = 14000
$ws.Cells.Item(72, 12).Value = 15000
$ws.Cells.Item(72, 13).Value = 14500
$ws.Cells.Item(72, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(72, 16).Value = 580
$ws.Cells.Item(73, 4).Value = 45085
$ws.Cells.Item(73, 10).Value = 700
$ws.Cells.Item(73, 11).Value = 9000
$ws.Cells.Item(73, 12).Value = 10000
$ws.Cells.Item(73, 13).Value = 9500
$ws.Cells.Item(73, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(73, 16).Value = 380
$ws.Cells.Item(74, 4).Value = 44825
$ws.Cells.Item(74, 10).Value = 700
$ws.Cells.Item(74, 11).Value = 8000
$ws.Cells.Item(74, 12).Value = 9000
$ws.Cells.Item(74, 13).Value = 8500
$ws.Cells.Item(74, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(74, 16).Value = 340
$ws.Cells.Item(75, 4).Value = 44504
$ws.Cells.Item(75, 10).Value = 700
$ws.Cells.Item(75, 11).Value = 6000
$ws.Cells.Item(75, 12).Value = 7000
$ws.Cells.Item(75, 13).Value = 6500
$ws.Cells.Item(75, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(75, 16).Value = 260
$ws.Cells.Item(76, 4).Value = 44742
$ws.Cells.Item(76, 10).Value = 600
$ws.Cells.Item(76, 11).Value = 13000
$ws.Cells.Item(76, 12).Value = 15000
$ws.Cells.Item(76, 13).Value = 14000
$ws.Cells.Item(76, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(76, 16).Value = 560
$ws.Cells.Item(77, 4).Value = 44490
$ws.Cells.Item(77, 10).Value = 400
$ws.Cells.Item(77, 11).Value = 5000
$ws.Cells.Item(77, 12).Value = 6000
$ws.Cells.Item(77, 13).Value = 5500
$ws.Cells.Item(77, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(77, 16).Value = 220
$ws.Cells.Item(78, 4).Value = 45071
$ws.Cells.Item(78, 10).Value = 500
$ws.Cells.Item(78, 11).Value = 10000
$ws.Cells.Item(78, 12).Value = 11000
$ws.Cells.Item(78, 13).Value = 10500
$ws.Cells.Item(78, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(78, 16).Value = 420
$ws.Cells.Item(79, 4).Value = 44902
$ws.Cells.Item(79, 10).Value = 500
$ws.Cells.Item(79, 11).Value = 5000
$ws.Cells.Item(79, 12).Value = 6000
$ws.Cells.Item(79, 13).Value = 5500
$ws.Cells.Item(79, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(79, 16).Value = 220
$ws.Cells.Item(80, 4).Value = 45070
$ws.Cells.Item(80, 10).Value = 1700
$ws.Cells.Item(80, 11).Value = 10000
$ws.Cells.Item(80, 12).Value = 11000
$ws.Cells.Item(80, 13).Value = 10500
$ws.Cells.Item(80, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(80, 16).Value = 420
$ws.Cells.Item(81, 4).Value = 44790
$ws.Cells.Item(81, 10).Value = 1460
$ws.Cells.Item(81, 11).Value = 9000
$ws.Cells.Item(81, 12).Value = 10000
$ws.Cells.Item(81, 13).Value = 9500
$ws.Cells.Item(81, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(81, 16).Value = 380
$ws.Cells.Item(82, 4).Value = 44714
$ws.Cells.Item(82, 10).Value = 500
$ws.Cells.Item(82, 11).Value = 12000
$ws.Cells.Item(82, 12).Value = 13000
$ws.Cells.Item(82, 13).Value = 12500
$ws.Cells.Item(82, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(82, 16).Value = 500
$ws.Cells.Item(83, 4).Value = 44826
$ws.Cells.Item(83, 10).Value = 1200
$ws.Cells.Item(83, 11).Value = 7000
$ws.Cells.Item(83, 12).Value = 8000
$ws.Cells.Item(83, 13).Value = 7500
$ws.Cells.Item(83, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(83, 16).Value = 300
$ws.Cells.Item(84, 4).Value = 44692
$ws.Cells.Item(84, 10).Value = 500
$ws.Cells.Item(84, 11).Value = 15000
$ws.Cells.Item(84, 12).Value = 16000
$ws.Cells.Item(84, 13).Value = 15500
$ws.Cells.Item(84, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(84, 16).Value = 620
$ws.Cells.Item(85, 4).Value = 44468
$ws.Cells.Item(85, 10).Value = 700
$ws.Cells.Item(85, 11).Value = 5000
$ws.Cells.Item(85, 12).Value = 6000
$ws.Cells.Item(85, 13).Value = 5500
$ws.Cells.Item(85, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(85, 16).Value = 220
$ws.Cells.Item(86, 4).Value = 45106
$ws.Cells.Item(86, 10).Value = 1100
$ws.Cells.Item(86, 11).Value = 11000
$ws.Cells.Item(86, 12).Value = 12000
$ws.Cells.Item(86, 13).Value = 11500
$ws.Cells.Item(86, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(86, 16).Value = 460
$ws.Cells.Item(87, 4).Value = 44391
$ws.Cells.Item(87, 10).Value = 500
$ws.Cells.Item(87, 11).Value = 9000
$ws.Cells.Item(87, 12).Value = 10000
$ws.Cells.Item(87, 13).Value = 9500
$ws.Cells.Item(87, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(87, 16).Value = 380
$ws.Cells.Item(88, 4).Value = 44420
$ws.Cells.Item(88, 10).Value = 1000
$ws.Cells.Item(88, 11).Value = 10000
$ws.Cells.Item(88, 12).Value = 11000
$ws.Cells.Item(88, 13).Value = 10500
$ws.Cells.Item(88, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(88, 16).Value = 420
$ws.Cells.Item(89, 4).Value = 44364
$ws.Cells.Item(89, 10).Value = 700
$ws.Cells.Item(89, 11).Value = 11000
$ws.Cells.Item(89, 12).Value = 12000
$ws.Cells.Item(89, 13).Value = 11500
$ws.Cells.Item(89, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(89, 16).Value = 460
$ws.Cells.Item(90, 4).Value = 44811
$ws.Cells.Item(90, 10).Value = 500
$ws.Cells.Item(90, 11).Value = 5000
$ws.Cells.Item(90, 12).Value = 6000
$ws.Cells.Item(90, 13).Value = 5500
$ws.Cells.Item(90, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(90, 16).Value = 220
$ws.Cells.Item(91, 4).Value = 45112
$ws.Cells.Item(91, 10).Value = 1200
$ws.Cells.Item(91, 11).Value = 12000
$ws.Cells.Item(91, 12).Value = 13000
$ws.Cells.Item(91, 13).Value = 12500
$ws.Cells.Item(91, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(91, 16).Value = 500
$ws.Cells.Item(92, 4).Value = 44749
$ws.Cells.Item(92, 10).Value = 1100
$ws.Cells.Item(92, 11).Value = 12000
$ws.Cells.Item(92, 12).Value = 14000
$ws.Cells.Item(92, 13).Value = 13000
$ws.Cells.Item(92, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(92, 16).Value = 520
$ws.Cells.Item(93, 4).Value = 44763
$ws.Cells.Item(93, 10).Value = 1100
$ws.Cells.Item(93, 11).Value = 11000
$ws.Cells.Item(93, 12).Value = 12000
$ws.Cells.Item(93, 13).Value = 11500
$ws.Cells.Item(93, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(93, 16).Value = 460

# Append two new rows (94, 95) with the same constant columns as the rest of the dataset
$ws.Cells.Item(94, 1).Value = 2
$ws.Cells.Item(94, 2).Value = 'Comercializadora del Agro de Limarí'
$ws.Cells.Item(94, 3).Value = 'Coquimbo'
$ws.Cells.Item(94, 4).Value = 44741
$ws.Cells.Item(94, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(94, 5).Value = 4
$ws.Cells.Item(94, 6).Value = 100112026
$ws.Cells.Item(94, 7).Value = 'Haba'
$ws.Cells.Item(94, 8).Value = 'Sin especificar'
$ws.Cells.Item(94, 9).Value = 'Primera'
$ws.Cells.Item(94, 10).Value = 700
$ws.Cells.Item(94, 11).Value = 14000
$ws.Cells.Item(94, 12).Value = 15000
$ws.Cells.Item(94, 13).Value = 14500
$ws.Cells.Item(94, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(94, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(94, 16).Value = 580
$ws.Cells.Item(94, 17).Value = 25
$ws.Cells.Item(94, 18).Value = 'Hortaliza'

$ws.Cells.Item(95, 1).Value = 2
$ws.Cells.Item(95, 2).Value = 'Comercializadora del Agro de Limarí'
$ws.Cells.Item(95, 3).Value = 'Coquimbo'
$ws.Cells.Item(95, 4).Value = 44777
$ws.Cells.Item(95, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(95, 5).Value = 4
$ws.Cells.Item(95, 6).Value = 100112026
$ws.Cells.Item(95, 7).Value = 'Haba'
$ws.Cells.Item(95, 8).Value = 'Sin especificar'
$ws.Cells.Item(95, 9).Value = 'Primera'
$ws.Cells.Item(95, 10).Value = 600
$ws.Cells.Item(95, 11).Value = 9000
$ws.Cells.Item(95, 12).Value = 10000
$ws.Cells.Item(95, 13).Value = 9500
$ws.Cells.Item(95, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(95, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(95, 16).Value = 380
$ws.Cells.Item(95, 17).Value = 25
$ws.Cells.Item(95, 18).Value = 'Hortaliza'

